$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: SCD0338 -> SCD0026
$ws.Name = "SCD0026"

# Update the TC_ID value (column B) on the two data rows from SCD0338-013 to SCD0026-013
$ws.Range("B2").Value = "SCD0026-013"
$ws.Range("B3").Value = "SCD0026-013"

# Reset the view: scroll back so column A is the left-most visible column,
# and move the active selection from Q4 to B4.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("B4").Select()
